$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 7 of this placeholder holds "27/10/2017 – UAT".
# Split the single run into three runs matching the new wording:
#   "27/10/2017 – " / "Teste de " / "Aceitação"
$para7 = $tr.Paragraphs(7, 1)
$para7.Text = "27/10/2017 – "
$run2 = $para7.InsertAfter("Teste de ")
$run3 = $run2.InsertAfter("Aceitação")
